# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# sheets, which hold identical data in this workbook.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 155
    "F3"  = 7081
    "F4"  = 4440
    "F10" = 66
    "F11" = 62
    "F14" = 117
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
